# Updated after RI IATI tag added by FCDO
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IATI activity IDs")

# iati_id values (column A) of the rows that need to be removed - these
# activities no longer belong in the partner activities list.
$idsToRemove = @(
    "GB-COH-RC000797-GB-GOV-1-300484-RIAH",
    "US-EIN-33-1112770-DFID_192010_LEIA_Business_Case",
    "GB-CHC-228248-F0192300",
    "XM-DAC-47021-1394-DFID",
    "GB-CHC-222655-PO6407",
    "GB-CHC-222655-PO5247REBUILD",
    "GB-CHC-1177110-R2HC",
    "CA-CRA_ACR-101182509-DFID-REGIONAL-SPARC",
    "US-EIN-91-1157127-PATH CENTRE OF EXCELLENCE FOR MICROARRAY PATCH PLATFORM DEVELOPMENT (MAPS) (300341-112)",
    "US-EIN-91-1157127-DAWN"
)

# Find the last used row in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 1 }

# Walk bottom-to-top so deleting a row doesn't shift the ones we still
# need to inspect.
for ($r = $lastRow; $r -ge 2; $r--) {
    $val = $ws.Cells.Item($r, 1).Value()
    if ($idsToRemove -contains $val) {
        $ws.Rows.Item($r).Delete()
    }
}

# Append the new activity row (FCDO fully funded tag added by FCDO for
# this new RI-tagged activity).
$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
$ws.Cells.Item($newRow, 1).Value = "GB-CHC-209131-A05500"
$ws.Cells.Item($newRow, 2).Value = "Foreign, Commonwealth and Development Office"
$ws.Cells.Item($newRow, 3).Value = "GB-1-204043"
$ws.Cells.Item($newRow, 4).Value = "FCDO fully funded"
$ws.Cells.Item($newRow, 5).Value = "British Council"

# Match the saved selection state (active cell D1).
$ws.Range("D1").Select()
